$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B9").Value = 6.345999999999999
$ws.Range("B18").Value = 5.126
$ws.Range("B20").Value = 6.24
$ws.Range("B27").Value = 6.161
$ws.Range("B35").Value = 8.317
$ws.Range("B69").Value = 5.306999999999999
$ws.Range("B76").Value = 6.308
$ws.Range("B78").Value = 7.811
$ws.Range("B82").Value = 5.366000000000001
$ws.Range("B83").Value = 5.667
$ws.Range("B93").Value = 5.610999999999999
